# Fill in the "Design" row (row 3) of the schedule with the same
# Expected/Actual Start/End dates style as the "SRS- Release" row (row 2),
# one day later (2022-09-08 instead of 2022-09-07), and move the active
# selection to C4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date-cell formatting from row 2 (C2:F2) onto row 3 (C3:F3) so
# the new values inherit the same number format / font / border as the
# rest of the date columns.
$ws.Range("C2:F2").Copy()
$ws.Range("C3:F3").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Write the actual date values (serial 44812 == 2022-09-08).
$ws.Range("C3").Value = 44812
$ws.Range("D3").Value = 44812
$ws.Range("E3").Value = 44812
$ws.Range("F3").Value = 44812

# Row now needs to grow slightly to fit the new content, matching row 2.
$ws.Rows("3").RowHeight = 14.25

# Move the selection/active cell.
$ws.Range("C4").Select()
